# Commit: "added calc for number of boards"
#
# Inserts two new columns (Min Order / James Ordered) after the existing
# Quantity-related columns, adds a "Boards" multiplier cell, and fills in a
# ROUNDUP-based "Min Order" quantity formula (with 5% buffer) down the BOM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at G:H - this shifts the old G (Supplied),
# H (Farnell Number), I (Digikey) and J (Comments) columns right to
# I, J, K and L respectively, carrying their values/styles/formats along.
$ws.Range("G:H").Insert()

# New header labels for the inserted columns (row 4).
# Order matters here: it controls the order the new strings are appended
# to the shared string table (Min Order, Boards, James Ordered).
$ws.Range("G4").Value = "Min Order"
$ws.Range("F3").Value = "Boards"
$ws.Range("H4").Value = "James Ordered"

# Number of boards being built - used as the multiplier in the Min Order
# calculation below.
$ws.Range("G3").Value = 5

# Min Order quantity = ROUNDUP(1.05 * (Quantity per board * number of boards), 0)
$ws.Range("G6").Formula = "=ROUNDUP(1.05*(C6*`$G`$3),0)"
$ws.Range("G7:G28").Formula = "=ROUNDUP(1.05*(C7*`$G`$3),0)"

# Recalculate so the cached formula results stored in the file are correct.
$excel.CalculateFullRebuild()
$excel.Calculate()

# Match the author's last active selection.
[void]$ws.Range("E19").Select()
